$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells stay as text (avoid Excel auto-numeric conversion)
$priceCells = @('D2', 'D3', 'D4', 'D5', 'D6', 'D8', 'D13', 'D14', 'D16', 'D17', 'D18', 'D19', 'D20', 'D22', 'D23', 'D24', 'D25', 'D27', 'D28', 'D29', 'D30', 'D31', 'D32', 'D34', 'D37', 'D38', 'D39', 'D40', 'D41', 'D42', 'D43', 'D44', 'D45', 'D46', 'D47', 'D48', 'D49')
foreach ($ref in $priceCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply updated values
$ws.Range('D2').Value = '60.534.50'
$ws.Range('E2').Value = '  +2.86%  '
$ws.Range('D3').Value = '3.030.93'
$ws.Range('E3').Value = '  +1.20%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.39%  '
$ws.Range('D5').Value = '571.53'
$ws.Range('E5').Value = '  +1.85%  '
$ws.Range('D6').Value = '141.26'
$ws.Range('E6').Value = '  +2.67%  '
$ws.Range('E7').Value = '  -0.18%  '
$ws.Range('D8').Value = '3.030.06'
$ws.Range('E8').Value = '  +1.58%  '
$ws.Range('E9').Value = '  +1.14%  '
$ws.Range('E10').Value = '  +4.17%  '
$ws.Range('E11').Value = '  +11.25%  '
$ws.Range('E12').Value = '  +0.38%  '
$ws.Range('D13').Value = '0.0000236'
$ws.Range('E13').Value = '  +2.50%  '
$ws.Range('D14').Value = '34.61'
$ws.Range('E14').Value = '  +2.59%  '
$ws.Range('E15').Value = '  -0.30%  '
$ws.Range('D16').Value = '3.544.49'
$ws.Range('E16').Value = '  +1.59%  '
$ws.Range('D17').Value = '7.16'
$ws.Range('E17').Value = '  +2.27%  '
$ws.Range('D18').Value = '3.031.89'
$ws.Range('E18').Value = '  +1.24%  '
$ws.Range('D19').Value = '60.530.60'
$ws.Range('E19').Value = '  +2.69%  '
$ws.Range('D20').Value = '444.07'
$ws.Range('E20').Value = '  +4.00%  '
$ws.Range('E21').Value = '  +1.26%  '
$ws.Range('D22').Value = '0.729'
$ws.Range('E22').Value = '  +2.28%  '
$ws.Range('D23').Value = '7.19'
$ws.Range('E23').Value = '  +0.69%  '
$ws.Range('D24').Value = '13.44'
$ws.Range('E24').Value = '  -0.06%  '
$ws.Range('D25').Value = '81.34'
$ws.Range('E25').Value = '  +1.20%  '
$ws.Range('E26').Value = '  +0.07%  '
$ws.Range('D27').Value = '2.26'
$ws.Range('E27').Value = '  +7.71%  '
$ws.Range('D28').Value = '0.997'
$ws.Range('E28').Value = '  -0.52%  '
$ws.Range('B29').Value = 'PancakeSwap'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D29').Value = '2.59'
$ws.Range('E29').Value = '  +2.45%  '
$ws.Range('B30').Value = 'RenderToken'
$ws.Range('C30').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D30').Value = '7.98'
$ws.Range('E30').Value = '  +2.97%  '
$ws.Range('B31').Value = 'EthereumClassic'
$ws.Range('C31').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D31').Value = '26.30'
$ws.Range('E31').Value = '  +2.23%  '
$ws.Range('B32').Value = 'NEARProtocol'
$ws.Range('C32').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D32').Value = '6.35'
$ws.Range('E32').Value = '  +3.91%  '
$ws.Range('E33').Value = '  +7.45%  '
$ws.Range('D34').Value = '0.0₃0804'
$ws.Range('E34').Value = '  +7.67%  '
$ws.Range('E35').Value = '  +4.40%  '
$ws.Range('E36').Value = '  +4.68%  '
$ws.Range('D37').Value = '2.12'
$ws.Range('E37').Value = '  +1.59%  '
$ws.Range('D38').Value = '49.47'
$ws.Range('E38').Value = '  +1.29%  '
$ws.Range('D39').Value = '2.95'
$ws.Range('E39').Value = '  +7.68%  '
$ws.Range('D40').Value = '8.73'
$ws.Range('E40').Value = '  -1.39%  '
$ws.Range('D41').Value = '405.73'
$ws.Range('E41').Value = '  +2.64%  '
$ws.Range('D42').Value = '0.0358'
$ws.Range('E42').Value = '  +2.63%  '
$ws.Range('D43').Value = '2.769.44'
$ws.Range('E43').Value = '  +1.55%  '
$ws.Range('D44').Value = '0.106'
$ws.Range('E44').Value = '  -1.35%  '
$ws.Range('D45').Value = '0.258'
$ws.Range('E45').Value = '  +5.09%  '
$ws.Range('D46').Value = '36.39'
$ws.Range('E46').Value = '  +13.55%  '
$ws.Range('D47').Value = '0.998'
$ws.Range('E47').Value = '  -0.08%  '
$ws.Range('D48').Value = '2.05'
$ws.Range('E48').Value = '  +1.40%  '
$ws.Range('D49').Value = '122.83'
$ws.Range('E49').Value = '  -2.34%  '
$ws.Range('E50').Value = '  +0.89%  '
$ws.Range('E51').Value = '  +1.25%  '
